$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number to Excel need to be
# forced to stay text (matching the original inlineStr string cells).
# Strategy: temporarily mark the cell as Text ("@") so the assignment is
# not reinterpreted as a number, then restore the original (Normal)
# style so no visible/structural formatting change is left behind.
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.989.35"
$ws.Range("E2").Value = "  +1.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.848.85"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.52%  "

# Row 5 - USDC
$ws.Range("D5").Value = "1.013"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6 - BNB
$ws.Range("D6").Value = "310.06"
$ws.Range("E6").Value = "  +0.28%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4783"
$ws.Range("E7").Value = "  +2.19%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3679"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07234"
$ws.Range("E9").Value = "  +1.31%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.9284"
$ws.Range("E10").Value = "  +2.80%  "

# Row 11 - Solana
$ws.Range("D11").Value = "19.71"
$ws.Range("E11").Value = "  +1.43%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.07742"
$ws.Range("E12").Value = "  +0.88%  "

# Row 13 and 14 swap: Polkadot <-> WrappedEther, with updated values
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.792.75"
$ws.Range("E13").Value = "  -3.00%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.348"
$ws.Range("E14").Value = "  +1.68%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "6.436"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "88.79"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "1.015"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000008632"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.46%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "27.010.19"
$ws.Range("E20").Value = "  +1.37%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.064"
$ws.Range("E22").Value = "  +0.84%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  +1.00%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "1.930"
$ws.Range("E24").Value = "  +0.85%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.92"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +1.68%  "

# Row 27 - LidoDAOToken
$ws.Range("D27").Value = "1.999"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28 - BitcoinCash
$ws.Range("D28").Value = "114.28"
$ws.Range("E28").Value = "  +0.51%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "4.959"
$ws.Range("E29").Value = "  +1.72%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.08888"
$ws.Range("E30").Value = "  +0.82%  "

# Row 31 - HuobiToken
$ws.Range("D31").Value = "3.328"
$ws.Range("E31").Value = "  +5.54%  "

# Row 32 - ARBITRUM
$ws.Range("D32").Value = "1.171"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "0.7414"
$ws.Range("E33").Value = "  +0.86%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "4.501"
$ws.Range("E34").Value = "  +1.56%  "

# Row 35 - RenderToken
$ws.Range("D35").Value = "2.737"
$ws.Range("E35").Value = "  -3.84%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "1.115"
$ws.Range("E36").Value = "  +3.84%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +1.55%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "0.05277"
$ws.Range("E38").Value = "  +2.41%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.979"

# Row 40 - TheSandbox
$ws.Range("D40").Value = "0.5212"
$ws.Range("E40").Value = "  +2.93%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "6.987"
$ws.Range("E41").Value = "  +1.74%  "

# Row 42 - Algorand
$ws.Range("D42").Value = "0.1512"
$ws.Range("E42").Value = "  +0.91%  "

# Row 43 - Aptos
$ws.Range("D43").Value = "8.229"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "10.60"
$ws.Range("E44").Value = "  +5.42%  "

# Row 45 - Decentraland
$ws.Range("E45").Value = "  +2.07%  "

# Row 46 - PaxDollar
$ws.Range("D46").Value = "1.014"
$ws.Range("E46").Value = "  +0.55%  "

# Row 47 - Quant
$ws.Range("D47").Value = "101.63"
$ws.Range("E47").Value = "  +3.15%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  +2.38%  "

# Row 49 - Aave
$ws.Range("D49").Value = "65.38"
$ws.Range("E49").Value = "  +2.21%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.06067"
$ws.Range("E50").Value = "  +0.72%  "

# Row 51 - EOS
$ws.Range("D51").Value = "0.8883"
$ws.Range("E51").Value = "  +4.24%  "

# Restore the Normal style on the text-forced cells so no stray
# number-format / style change is left on the cell.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
